# Auto-generated edit script: apply cryptos.xlsx value updates via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.145.08"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.303.92"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.86%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.64%  "
$ws.Range("E7").Value = "  -4.87%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.296.13"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.611"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.70"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.57%  "
$ws.Range("E13").Value = "  -6.01%  "
$ws.Range("E14").Value = "  -6.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.835.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.86%  "
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.325.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.71%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.50%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.139.03"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.959"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "402.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.50"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.65%  "
$ws.Range("E28").Value = "  -7.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.84"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "570.74"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.02%  "
$ws.Range("E34").Value = "  -6.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -3.64%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.61%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.134.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0728"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -10.44%  "
$ws.Range("E42").Value = "  -6.54%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0398"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.80%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.07%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.42"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.43%  "
